$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 will hold the new leave application record (same data as the old
# row 4 / Animesh Roy, but with a new Leave ID and refreshed dates).
# Force text formatting first so Excel does not auto-convert the
# date-like strings ("02-12-2025", "2025-07-29", ...) into date serials.
$ws.Range("A2:D2").NumberFormat = "@"
$ws.Range("F2:M2").NumberFormat = "@"
$ws.Range("N2:O2").NumberFormat = "@"

$ws.Range("A2").Value = "LV-1753795678579-DDHA"
$ws.Range("B2").Value = "Manager"
$ws.Range("C2").Value = "PILLP305"
$ws.Range("D2").Value = "Animesh Roy"
$ws.Range("E2").Value = 8145312848
$ws.Range("F2").Value = "animesh.roy@pillp.in"
$ws.Range("G2").Value = "Mechanical HOD"
$ws.Range("H2").Value = "Mechanical"
$ws.Range("I2").Value = "Raichur"
$ws.Range("J2").Value = "02-12-2025"
$ws.Range("K2").Value = "07-12-2025"
$ws.Range("L2").Value = "For my personal reason"
$ws.Range("M2").Value = "Pending"
$ws.Range("N2").Value = "2025-07-29"
$ws.Range("O2").Value = "2025-07-29"

# Restore the default (Normal) cell style so the text-format override
# used above doesn't leave a lingering number format on the cells.
$ws.Range("A2:D2").Style = "Normal"
$ws.Range("F2:M2").Style = "Normal"
$ws.Range("N2:O2").Style = "Normal"

# The old rows 3 and 4 (which held the other two leave applications) are
# removed entirely; remaining rows shift up so the sheet now only spans
# down to row 2.
$ws.Rows("3:4").Delete()
